$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42606.882638888892
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"

$ws.Range("B4").Value = 28
$ws.Range("C4").Value = 66
$ws.Range("D4").Value = 33
$ws.Range("E4").Value = 66
$ws.Range("F4").Value = 33
$ws.Range("G4").Value = 7587
$ws.Range("H4").Value = 11117
$ws.Range("I4").Value = 2095
$ws.Range("J4").Value = 221
$ws.Range("K4").Value = 112
$ws.Range("L4").Value = 12
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = "Noun"
